$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 94.5
$ws.Range("I11").Value = 94.5
$ws.Range("K11").Value = 94.5
$ws.Range("M11").Value = 45.5
$ws.Range("H34").Value = 5591.857
$ws.Range("I34").Value = 2323.5
$ws.Range("J34").Value = 9949.666999999999
$ws.Range("K34").Value = 2323.5
$ws.Range("L34").Value = 9949.666999999999
$ws.Range("M34").Value = -2120.5
$ws.Range("N34").Value = -10355.667
$ws.Range("H36").Value = 5591.857
$ws.Range("I36").Value = 2323.5
$ws.Range("J36").Value = 9949.666999999999
$ws.Range("K36").Value = 2323.5
$ws.Range("L36").Value = 9949.666999999999
$ws.Range("M36").Value = -1608.5
$ws.Range("N36").Value = -11379.667
$ws.Range("H100").Value = 1686.4286
$ws.Range("I100").Value = 1761
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1761
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1220
$ws.Range("N100").Value = -2582
$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880
$ws.Range("H129").Value = 1076.3773
$ws.Range("J129").Value = 1140.7916
$ws.Range("L129").Value = 3422.3748
$ws.Range("N129").Value = -13422.3748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7513.8784
$ws.Range("I32").Value = 5205.5
$ws.Range("K32").Value = 5205.5
$ws.Range("M32").Value = -4918.5
$ws.Range("H61").Value = 252911.78
$ws.Range("I61").Value = 2158
$ws.Range("J61").Value = 530060.7
$ws.Range("K61").Value = 2158
$ws.Range("L61").Value = 530060.7
$ws.Range("M61").Value = -1946
$ws.Range("N61").Value = -530484.7
$ws.Range("H110").Value = 1531.65
$ws.Range("I110").Value = 1310.2858
$ws.Range("K110").Value = 1310.2858
$ws.Range("M110").Value = 734.7141999999999
$ws.Range("H136").Value = 252911.78
$ws.Range("I136").Value = 2158
$ws.Range("J136").Value = 530060.7
$ws.Range("K136").Value = 6474
$ws.Range("L136").Value = 1590182.1
$ws.Range("M136").Value = -3924
$ws.Range("N136").Value = -1595282.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 19071.084
$ws.Range("I105").Value = 22291.5
$ws.Range("K105").Value = 22291.5
$ws.Range("M105").Value = -20544.5
$ws.Range("H134").Value = 2267.8823
$ws.Range("I134").Value = 1846.4
$ws.Range("K134").Value = 5539.200000000001
$ws.Range("M134").Value = -3004.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1668.28
$ws.Range("I16").Value = 1310.0555
$ws.Range("K16").Value = 1310.0555
$ws.Range("M16").Value = -1023.0555
$ws.Range("H31").Value = 3369.0176
$ws.Range("I31").Value = 2326.261
$ws.Range("J31").Value = 4074.4119
$ws.Range("K31").Value = 2326.261
$ws.Range("L31").Value = 4074.4119
$ws.Range("M31").Value = -2031.261
$ws.Range("N31").Value = -4664.4119
$ws.Range("H34").Value = 3369.0176
$ws.Range("I34").Value = 2326.261
$ws.Range("J34").Value = 4074.4119
$ws.Range("K34").Value = 2326.261
$ws.Range("L34").Value = 4074.4119
$ws.Range("M34").Value = -2124.261
$ws.Range("N34").Value = -4478.4119
$ws.Range("H113").Value = 1668.28
$ws.Range("I113").Value = 1310.0555
$ws.Range("K113").Value = 1310.0555
$ws.Range("M113").Value = 859.9445000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 656.13043
$ws.Range("J107").Value = 680.5
$ws.Range("L107").Value = 2041.5
$ws.Range("N107").Value = -5881.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 18884.166
$ws.Range("I80").Value = 51402.5
$ws.Range("J80").Value = 2625
$ws.Range("K80").Value = 51402.5
$ws.Range("L80").Value = 2625
$ws.Range("M80").Value = -50404.5
$ws.Range("N80").Value = -4621
$ws.Range("H83").Value = 18884.166
$ws.Range("I83").Value = 51402.5
$ws.Range("J83").Value = 2625
$ws.Range("K83").Value = 257012.5
$ws.Range("L83").Value = 13125
$ws.Range("M83").Value = -252020.5
$ws.Range("N83").Value = -23109
$ws.Range("H113").Value = 1584.9231
$ws.Range("I113").Value = 1080.5333
$ws.Range("K113").Value = 1080.5333
$ws.Range("M113").Value = 1089.4667
$ws.Range("H132").Value = 2973
$ws.Range("I132").Value = 3057.8
$ws.Range("K132").Value = 9173.400000000001
$ws.Range("M132").Value = -6643.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3000
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 3000
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3214
$ws.Range("H132").Value = 2843.348
$ws.Range("I132").Value = 2705.158
$ws.Range("J132").Value = 3499.75
$ws.Range("K132").Value = 8115.474
$ws.Range("L132").Value = 10499.25
$ws.Range("M132").Value = -5585.474
$ws.Range("N132").Value = -15559.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 11000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 11000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 11000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -11826
$ws.Range("H42").Value = 6996.6665
$ws.Range("I42").Value = 1660
$ws.Range("J42").Value = 12333.333
$ws.Range("K42").Value = 1660
$ws.Range("L42").Value = 12333.333
$ws.Range("M42").Value = -1282
$ws.Range("N42").Value = -13089.333
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15298
$ws.Range("H100").Value = 1000002
$ws.Range("I100").Value = 1000002
$ws.Range("K100").Value = 2000004
$ws.Range("M100").Value = -1999463
$ws.Range("H132").Value = 1293.7407
$ws.Range("I132").Value = 913.1795
$ws.Range("J132").Value = 2283.2
$ws.Range("K132").Value = 2739.5385
$ws.Range("L132").Value = 6849.599999999999
$ws.Range("M132").Value = -209.5384999999997
$ws.Range("N132").Value = -11909.6
$ws.Range("H136").Value = 2173.093
$ws.Range("I136").Value = 2397.3044
$ws.Range("J136").Value = 1915.25
$ws.Range("K136").Value = 7191.9132
$ws.Range("L136").Value = 5745.75
$ws.Range("M136").Value = -4641.9132
$ws.Range("N136").Value = -10845.75
